$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (e.g. "26.50" -> 26.5, "1.00" -> 1), which would not match the
# original inline-string (text) cell contents.
$textCells = @("D5","D6","D7","D8","D13","D17","D19","D20","D21","D22","D23","D26","D27","D33","D36","D37","D39","D42","D43","D44","D45","D47","D48","D49","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "75.350.61"
$ws.Range("E2").Value = "  +7.77%  "
$ws.Range("D3").Value = "2.667.95"
$ws.Range("E3").Value = "  +9.18%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "187.09"
$ws.Range("E5").Value = "  +12.30%  "
$ws.Range("D6").Value = "587.49"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  +3.99%  "
$ws.Range("E9").Value = "  +12.09%  "
$ws.Range("D10").Value = "2.670.90"
$ws.Range("E10").Value = "  +9.15%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  +6.64%  "
$ws.Range("D13").Value = "4.73"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.162.49"
$ws.Range("E14").Value = "  +9.00%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "74.922.10"
$ws.Range("E15").Value = "  +7.26%  "
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "26.50"
$ws.Range("E17").Value = "  +9.86%  "
$ws.Range("D18").Value = "2.677.26"
$ws.Range("E18").Value = "  +9.47%  "
$ws.Range("D19").Value = "9.18"
$ws.Range("E19").Value = "  +28.52%  "
$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  +9.75%  "
$ws.Range("D21").Value = "371.59"
$ws.Range("E21").Value = "  +9.03%  "
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  +13.89%  "
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  +4.95%  "
$ws.Range("E24").Value = "  +3.74%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "69.80"
$ws.Range("E26").Value = "  +5.17%  "
$ws.Range("D27").Value = "4.14"
$ws.Range("E27").Value = "  +8.33%  "
$ws.Range("E28").Value = "  +10.01%  "
$ws.Range("E29").Value = "  +8.43%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "0.0₃0945"
$ws.Range("E31").Value = "  +10.41%  "
$ws.Range("E32").Value = "  +14.35%  "
$ws.Range("D33").Value = "520.47"
$ws.Range("E33").Value = "  +12.94%  "
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("E35").Value = "  +7.90%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "163.52"
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("E38").Value = "  +5.88%  "
$ws.Range("D39").Value = "19.18"
$ws.Range("E39").Value = "  +5.26%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "169.99"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  +12.82%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  +10.16%  "
$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  +8.65%  "
$ws.Range("E46").Value = "  +9.23%  "
$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  +12.22%  "
$ws.Range("D48").Value = "39.09"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "0.0846"
$ws.Range("E49").Value = "  +16.63%  "
$ws.Range("E50").Value = "  +7.28%  "
$ws.Range("D51").Value = "0.534"
$ws.Range("E51").Value = "  +8.82%  "
